# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 13:29"

# India (row 6) - updated covid numbers
$ws.Range("B6").Value = 1588177
$ws.Range("C6").Value = 3793
$ws.Range("E6").Value = 530535

# Alemania (row 21)
$ws.Range("B21").Value = 208885
$ws.Range("C21").Value = 74
$ws.Range("E21").Value = 7673

# Catar (row 26)
$ws.Range("B26").Value = 110460
$ws.Range("C26").Value = 307
$ws.Range("D26").Value = 107135
$ws.Range("E26").Value = 3154
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 171

# Bielorrusia (row 38)
$ws.Range("B38").Value = 67665
$ws.Range("C38").Value = 147
$ws.Range("D38").Value = 61765
$ws.Range("E38").Value = 5347
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 553

# Belgica (row 39)
$ws.Range("D39").Value = 17492
$ws.Range("E39").Value = 40007

# Venezuela (row 70)
$ws.Range("E70").Value = 6583
$ws.Range("H70").Value = 154

# Row 83: now Madagascar, overtaking Republica de Macedonia with fresh numbers
$ws.Range("A83").Value = "Madagascar"
$ws.Range("B83").Value = 10748
$ws.Range("C83").Value = 431
$ws.Range("D83").Value = 7461
$ws.Range("E83").Value = 3182
$ws.Range("G83").Value = 6
$ws.Range("H83").Value = 105

# Row 84: now Republica de Macedonia, keeping its previous numbers
$ws.Range("A84").Value = "Republica de Macedonia"
$ws.Range("B84").Value = 10503
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 5931
$ws.Range("E84").Value = 4096
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 476

# Malta (row 155)
$ws.Range("B155").Value = 729
$ws.Range("C155").Value = 9
$ws.Range("E155").Value = 55

# Vietnam (row 164)
$ws.Range("B164").Value = 464
$ws.Range("C164").Value = 5
$ws.Range("E164").Value = 95
